$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "25.551.68"
$ws.Range("D3").Value = "1.666.11"
$ws.Range("E3").Value = "  +0.79%  "
Set-TextValue $ws.Range("D4") "0.9998"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "236.79"
$ws.Range("E5").Value = "  +0.06%  "
Set-TextValue $ws.Range("D7") "0.4797"
$ws.Range("E7").Value = "  +0.28%  "
Set-TextValue $ws.Range("D8") "0.2632"
$ws.Range("E8").Value = "  +0.42%  "
Set-TextValue $ws.Range("D9") "0.06159"
$ws.Range("E9").Value = "  +3.03%  "
Set-TextValue $ws.Range("D10") "0.07085"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "1.667.32"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +3.08%  "
Set-TextValue $ws.Range("D13") "0.5963"
$ws.Range("E13").Value = "  -3.78%  "
Set-TextValue $ws.Range("D14") "4.399"
$ws.Range("E14").Value = "  -4.41%  "
Set-TextValue $ws.Range("D15") "74.52"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "25.555.67"
$ws.Range("E18").Value = "  +1.92%  "
Set-TextValue $ws.Range("D19") "0.000006769"
$ws.Range("E19").Value = "  +3.34%  "
Set-TextValue $ws.Range("D20") "11.45"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "1.880.58"
$ws.Range("E21").Value = "  +0.82%  "
Set-TextValue $ws.Range("D22") "4.461"
$ws.Range("E22").Value = "  +0.34%  "
Set-TextValue $ws.Range("D23") "8.690"
$ws.Range("E23").Value = "  +2.61%  "
Set-TextValue $ws.Range("D24") "5.338"
$ws.Range("E24").Value = "  +1.33%  "
Set-TextValue $ws.Range("D25") "134.93"
$ws.Range("E25").Value = "  +1.26%  "
Set-TextValue $ws.Range("D26") "15.10"
$ws.Range("E26").Value = "  +2.60%  "
Set-TextValue $ws.Range("D27") "1.404"
$ws.Range("E27").Value = "  +0.37%  "
Set-TextValue $ws.Range("D28") "104.91"
$ws.Range("E28").Value = "  +3.26%  "
Set-TextValue $ws.Range("D29") "1.688"
$ws.Range("E29").Value = "  -0.47%  "
Set-TextValue $ws.Range("D30") "3.978"
$ws.Range("E30").Value = "  +4.58%  "
Set-TextValue $ws.Range("D31") "3.664"
$ws.Range("E31").Value = "  +4.09%  "
Set-TextValue $ws.Range("D32") "0.07675"
$ws.Range("E32").Value = "  -2.90%  "
Set-TextValue $ws.Range("D33") "0.04336"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.47%  "
Set-TextValue $ws.Range("D36") "0.6153"
$ws.Range("E36").Value = "  +5.33%  "
Set-TextValue $ws.Range("D37") "0.9517"
$ws.Range("E37").Value = "  +1.16%  "
Set-TextValue $ws.Range("D38") "2.608"
$ws.Range("E38").Value = "  -0.50%  "
Set-TextValue $ws.Range("D39") "0.8621"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  -0.03%  "
Set-TextValue $ws.Range("D41") "0.01512"
$ws.Range("E41").Value = "  -1.56%  "
Set-TextValue $ws.Range("D42") "1.885"
$ws.Range("E42").Value = "  +3.06%  "
Set-TextValue $ws.Range("D43") "97.78"
$ws.Range("E43").Value = "  -0.71%  "
Set-TextValue $ws.Range("D44") "0.3774"
$ws.Range("E44").Value = "  +2.04%  "
Set-TextValue $ws.Range("D45") "4.690"
$ws.Range("E45").Value = "  -2.91%  "
Set-TextValue $ws.Range("D46") "0.1122"
$ws.Range("E46").Value = "  +0.45%  "
Set-TextValue $ws.Range("D47") "6.216"
$ws.Range("E47").Value = "  +2.75%  "
Set-TextValue $ws.Range("D48") "0.05264"
$ws.Range("E48").Value = "  +2.21%  "
Set-TextValue $ws.Range("D49") "29.52"
$ws.Range("E49").Value = "  +0.71%  "
Set-TextValue $ws.Range("D50") "7.395"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("E51").Value = "  +0.13%  "
